$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host $ws.Name
$v = $ws.Range("I8").Value
Write-Host $v
Write-Host "A11 sel"
Write-Host $ws.Range("A4").Value
